$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2131.946
$ws.Range("I15").Value = 2131.946
$ws.Range("K15").Value = 6395.838
$ws.Range("M15").Value = -6226.838
$ws.Range("H40").Value = 3999.8572
$ws.Range("J40").Value = 3999.926
$ws.Range("L40").Value = 3999.926
$ws.Range("N40").Value = -4349.925999999999
$ws.Range("H51").Value = 3862.1924
$ws.Range("J51").Value = 3059.5557
$ws.Range("L51").Value = 3059.5557
$ws.Range("N51").Value = -4027.5557
$ws.Range("H137").Value = 2170359.8
$ws.Range("J137").Value = 4276409
$ws.Range("L137").Value = 12829227
$ws.Range("N137").Value = -12834327

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 69994
$ws.Range("J44").Value = 69994
$ws.Range("L44").Value = 69994
$ws.Range("N44").Value = -70970
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H61").Value = 3014.1035
$ws.Range("I61").Value = 2164.611
$ws.Range("K61").Value = 2164.611
$ws.Range("M61").Value = -1952.611
$ws.Range("H74").Value = 1598.6744
$ws.Range("I74").Value = 1485.1515
$ws.Range("K74").Value = 1485.1515
$ws.Range("M74").Value = -611.1514999999999
$ws.Range("H77").Value = 1598.6744
$ws.Range("I77").Value = 1485.1515
$ws.Range("K77").Value = 7425.7575
$ws.Range("M77").Value = -3057.7575
$ws.Range("H122").Value = 6799.9565
$ws.Range("I122").Value = 6300
$ws.Range("J122").Value = 6938.8335
$ws.Range("K122").Value = 18900
$ws.Range("L122").Value = 20816.5005
$ws.Range("M122").Value = -16450
$ws.Range("N122").Value = -25716.5005
$ws.Range("H132").Value = 2420.0264
$ws.Range("J132").Value = 3497.5
$ws.Range("L132").Value = 10492.5
$ws.Range("N132").Value = -15552.5
$ws.Range("H136").Value = 3014.1035
$ws.Range("I136").Value = 2164.611
$ws.Range("K136").Value = 6493.833
$ws.Range("M136").Value = -3943.833

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4536
$ws.Range("I20").Value = 3633
$ws.Range("J20").Value = 5077.8
$ws.Range("K20").Value = 3633
$ws.Range("L20").Value = 5077.8
$ws.Range("M20").Value = -3386
$ws.Range("N20").Value = -5571.8
$ws.Range("H96").Value = 54577.4
$ws.Range("J96").Value = 93944
$ws.Range("L96").Value = 93944
$ws.Range("N96").Value = -99436
$ws.Range("H134").Value = 1641.2
$ws.Range("I134").Value = 1724.1111
$ws.Range("K134").Value = 5172.3333
$ws.Range("M134").Value = -2637.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3895.7058
$ws.Range("I31").Value = 1585.2413
$ws.Range("J31").Value = 6941.3184
$ws.Range("K31").Value = 1585.2413
$ws.Range("L31").Value = 6941.3184
$ws.Range("M31").Value = -1290.2413
$ws.Range("N31").Value = -7531.3184
$ws.Range("H34").Value = 3895.7058
$ws.Range("I34").Value = 1585.2413
$ws.Range("J34").Value = 6941.3184
$ws.Range("K34").Value = 1585.2413
$ws.Range("L34").Value = 6941.3184
$ws.Range("M34").Value = -1383.2413
$ws.Range("N34").Value = -7345.3184
$ws.Range("H87").Value = 82567.39999999999
$ws.Range("J87").Value = 82567.39999999999
$ws.Range("L87").Value = 82567.39999999999
$ws.Range("N87").Value = -84939.39999999999
$ws.Range("H90").Value = 82567.39999999999
$ws.Range("J90").Value = 82567.39999999999
$ws.Range("L90").Value = 247702.2
$ws.Range("N90").Value = -259558.2
$ws.Range("H134").Value = 2496.5
$ws.Range("I134").Value = 2017.6666
$ws.Range("K134").Value = 6052.9998
$ws.Range("M134").Value = -3517.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 57.5
$ws.Range("I13").Value = 57.5
$ws.Range("K13").Value = 172.5
$ws.Range("M13").Value = -4.5
$ws.Range("H15").Value = 147
$ws.Range("I15").Value = 90.5
$ws.Range("J15").Value = 184.66667
$ws.Range("K15").Value = 271.5
$ws.Range("L15").Value = 554.00001
$ws.Range("M15").Value = -131.5
$ws.Range("N15").Value = -834.00001
$ws.Range("H26").Value = 150
$ws.Range("I26").Value = 150
$ws.Range("K26").Value = 450
$ws.Range("M26").Value = -162
$ws.Range("H68").Value = 817.6
$ws.Range("I68").Value = 664
$ws.Range("K68").Value = 1992
$ws.Range("M68").Value = -1181
$ws.Range("H71").Value = 817.6
$ws.Range("I71").Value = 664
$ws.Range("K71").Value = 5976
$ws.Range("M71").Value = -1920
$ws.Range("H98").Value = 1587.8889
$ws.Range("J98").Value = 866.3333
$ws.Range("L98").Value = 2598.9999
$ws.Range("N98").Value = -5594.9999
$ws.Range("H103").Value = 639.6
$ws.Range("J103").Value = 500
$ws.Range("L103").Value = 1500
$ws.Range("N103").Value = -3258
$ws.Range("H129").Value = 1916.2727
$ws.Range("I129").Value = 308.33334
$ws.Range("J129").Value = 2519.25
$ws.Range("K129").Value = 925.0000200000001
$ws.Range("L129").Value = 7557.75
$ws.Range("M129").Value = 4074.99998
$ws.Range("N129").Value = -17557.75
$ws.Range("H134").Value = 2298.2415
$ws.Range("I134").Value = 2298.2415
$ws.Range("K134").Value = 6894.7245
$ws.Range("M134").Value = -1824.7245
$ws.Range("H136").Value = 1708.5555
$ws.Range("I136").Value = 1708.5555
$ws.Range("K136").Value = 5125.666499999999
$ws.Range("M136").Value = -25.66649999999936

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 27407.334
$ws.Range("J15").Value = 27407.334
$ws.Range("L15").Value = 27407.334
$ws.Range("N15").Value = -27983.334
$ws.Range("H81").Value = 27407.334
$ws.Range("J81").Value = 27407.334
$ws.Range("L81").Value = 27407.334
$ws.Range("N81").Value = -29403.334
$ws.Range("H84").Value = 27407.334
$ws.Range("J84").Value = 27407.334
$ws.Range("L84").Value = 82222.00199999999
$ws.Range("N84").Value = -92206.00199999999
$ws.Range("H132").Value = 3045.4736
$ws.Range("I132").Value = 2483.2144
$ws.Range("J132").Value = 4619.8
$ws.Range("K132").Value = 7449.6432
$ws.Range("L132").Value = 13859.4
$ws.Range("M132").Value = -4919.6432
$ws.Range("N132").Value = -18919.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2988.25
$ws.Range("J40").Value = 2976.5
$ws.Range("L40").Value = 2976.5
$ws.Range("N40").Value = -3248.5
$ws.Range("H61").Value = 1291.3334
$ws.Range("I61").Value = 1291.3334
$ws.Range("K61").Value = 1291.3334
$ws.Range("M61").Value = -1089.3334
$ws.Range("H113").Value = 1291.3334
$ws.Range("I113").Value = 1291.3334
$ws.Range("K113").Value = 1291.3334
$ws.Range("M113").Value = 878.6666
$ws.Range("H136").Value = 6182.381
$ws.Range("I136").Value = 4474.222
$ws.Range("J136").Value = 7463.5
$ws.Range("K136").Value = 13422.666
$ws.Range("L136").Value = 22390.5
$ws.Range("M136").Value = -10872.666
$ws.Range("N136").Value = -27490.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3420.0557
$ws.Range("I81").Value = 1708.6666
$ws.Range("J81").Value = 5131.4443
$ws.Range("K81").Value = 3417.3332
$ws.Range("L81").Value = 10262.8886
$ws.Range("M81").Value = -2356.3332
$ws.Range("N81").Value = -12384.8886
$ws.Range("H84").Value = 3420.0557
$ws.Range("I84").Value = 1708.6666
$ws.Range("J84").Value = 5131.4443
$ws.Range("K84").Value = 17086.666
$ws.Range("L84").Value = 51314.443
$ws.Range("M84").Value = -11782.666
$ws.Range("N84").Value = -61922.443
$ws.Range("H123").Value = 60000
$ws.Range("J123").Value = 70000
$ws.Range("L123").Value = 70000
$ws.Range("N123").Value = -79800
$ws.Range("H125").Value = 58715
$ws.Range("J125").Value = 58715
$ws.Range("L125").Value = 58715
$ws.Range("N125").Value = -68555
$ws.Range("H132").Value = 2398.3333
$ws.Range("I132").Value = 2490.9443
$ws.Range("J132").Value = 2120.5
$ws.Range("K132").Value = 7472.8329
$ws.Range("L132").Value = 6361.5
$ws.Range("M132").Value = -4942.8329
$ws.Range("N132").Value = -11421.5
$ws.Range("H138").Value = 97288.664
$ws.Range("J138").Value = 97288.664
$ws.Range("L138").Value = 97288.664
$ws.Range("N138").Value = -107568.664
